$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new row at position 14 (shifts old rows 14-22 down to 15-23)
# and populate it with the new Pinus monticola (PIMO) flammability data.
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = "PIMO"
$ws.Cells.Item(14, 2).Value = "Pinus monticola"
$ws.Cells.Item(14, 3).Value = 75.14
$ws.Cells.Item(14, 4).Value = 90.29
$ws.Cells.Item(14, 5).Value = 391.43
$ws.Cells.Item(14, 6).Value = 82.78
$ws.Cells.Item(14, 7).Value = "Banwell and Varner unpub data"

# Update the saved selection to match the author's final cursor position.
[void]$ws.Range("F1").Select()
